# Update Footprint X (P), Footprint Y (Q) and Footprint Rot (R) values on the
# "BoM" worksheet of the pedalboard-hw BoM workbook, reflecting a refreshed
# pedalboard-hw positional build (coordinates/rotations regenerated).
#
# The target cells store their numeric-looking values as plain text (e.g.
# "136.7734", "180.0000") rather than real numbers, so trailing zeros are
# preserved. Excel's automation model otherwise converts numeric-looking
# strings to real numbers, so we force each touched cell to Text format
# before writing the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# row -> column -> new value (only columns that actually change are listed)
$updates = @{
    16 = @{ P = "136.7734"; Q = "73.5802";  R = "180.0000" }
    20 = @{ P = "65.6560";  Q = "89.0730";  R = "270.0000" }
    30 = @{ P = "63.0160";  Q = "87.5730";  R = "0.0000" }
    33 = @{ P = "75.2034";  Q = "86.7030";  R = "0.0000" }
    35 = @{ P = "75.2034";  Q = "88.6030";  R = "180.0000" }
    36 = @{ P = "63.5410";  Q = "89.9980";  R = "90.0000" }
    42 = @{ P = "32.3744";  Q = "78.3300" }
    46 = @{ P = "55.0870";  Q = "79.1174" }
    51 = @{ P = "72.4560";  Q = "91.7230";  R = "0.0000" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
